$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1639.3572
$ws.Range("I6").Value = 1879.3334
$ws.Range("K6").Value = 5638.0002
$ws.Range("M6").Value = -5526.0002
# Row 8
$ws.Range("H8").Value = 37.666668
$ws.Range("I8").Value = 37.666668
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 113.000004
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 25.999996
$ws.Range("N8").ClearContents()
# Row 39
$ws.Range("H39").Value = 397.85715
$ws.Range("I39").Value = 308.875
$ws.Range("J39").Value = 516.5
$ws.Range("K39").Value = 926.625
$ws.Range("L39").Value = 1549.5
$ws.Range("M39").Value = -630.625
$ws.Range("N39").Value = -2141.5
# Row 64
$ws.Range("H64").Value = 4687.4375
$ws.Range("I64").Value = 4571.2856
$ws.Range("J64").Value = 4777.778
$ws.Range("K64").Value = 4571.2856
$ws.Range("L64").Value = 4777.778
$ws.Range("M64").Value = -4323.2856
$ws.Range("N64").Value = -5273.778
# Row 67
$ws.Range("H67").Value = 4687.4375
$ws.Range("I67").Value = 4571.2856
$ws.Range("J67").Value = 4777.778
$ws.Range("K67").Value = 4571.2856
$ws.Range("L67").Value = 4777.778
$ws.Range("M67").Value = -3713.2856
$ws.Range("N67").Value = -6493.778
# Row 74
$ws.Range("H74").Value = 252000
$ws.Range("I74").Value = 4000
$ws.Range("K74").Value = 4000
$ws.Range("M74").Value = -3064
# Row 77
$ws.Range("H77").Value = 252000
$ws.Range("I77").Value = 4000
$ws.Range("K77").Value = 20000
$ws.Range("M77").Value = -15320
# Row 87
$ws.Range("H87").Value = 43000
$ws.Range("J87").Value = 43000
$ws.Range("L87").Value = 43000
$ws.Range("N87").Value = -45496
# Row 90
$ws.Range("H90").Value = 43000
$ws.Range("J90").Value = 43000
$ws.Range("L90").Value = 129000
$ws.Range("N90").Value = -141480
# Row 137
$ws.Range("H137").Value = 1863.84
$ws.Range("I137").Value = 1358.5454
$ws.Range("K137").Value = 4075.6362
$ws.Range("M137").Value = -1525.6362

$ws = $wb.Worksheets.Item("ARM")
# Row 11
$ws.Range("H11").Value = 10001
$ws.Range("I11").Value = 9999
$ws.Range("J11").Value = 10003
$ws.Range("K11").Value = 9999
$ws.Range("L11").Value = 10003
$ws.Range("M11").Value = -9855
$ws.Range("N11").Value = -10291
# Row 44
$ws.Range("H44").Value = 33333.332
$ws.Range("J44").Value = 33333.332
$ws.Range("L44").Value = 33333.332
$ws.Range("N44").Value = -34309.332
# Row 80
$ws.Range("H80").Value = 39998.125
$ws.Range("J80").Value = 39998.125
$ws.Range("L80").Value = 39998.125
$ws.Range("N80").Value = -41994.125
# Row 83
$ws.Range("H83").Value = 39998.125
$ws.Range("J83").Value = 39998.125
$ws.Range("L83").Value = 119994.375
$ws.Range("N83").Value = -129978.375
# Row 96
$ws.Range("H96").Value = 25000
$ws.Range("J96").Value = 25000
$ws.Range("L96").Value = 25000
$ws.Range("N96").Value = -30492

$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Range("H81").Value = 48486.25
$ws.Range("J81").Value = 48486.25
$ws.Range("L81").Value = 48486.25
$ws.Range("N81").Value = -50608.25
# Row 82
$ws.Range("H82").Value = 33037.24
$ws.Range("J82").Value = 41037.5
$ws.Range("L82").Value = 41037.5
$ws.Range("N82").Value = -41803.5
# Row 84
$ws.Range("H84").Value = 48486.25
$ws.Range("J84").Value = 48486.25
$ws.Range("L84").Value = 145458.75
$ws.Range("N84").Value = -156066.75
# Row 85
$ws.Range("H85").Value = 33037.24
$ws.Range("J85").Value = 41037.5
$ws.Range("L85").Value = 41037.5
$ws.Range("N85").Value = -43689.5
# Row 107
$ws.Range("H107").Value = 1873.4286
$ws.Range("I107").Value = 1922.8
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 1922.8
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = -2.799999999999955
$ws.Range("N107").Value = -5590

$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 15921.444
$ws.Range("J41").Value = 22499
$ws.Range("L41").Value = 22499
$ws.Range("N41").Value = -23355
# Row 50
$ws.Range("H50").Value = 28512.834
$ws.Range("J50").Value = 29998.8
$ws.Range("L50").Value = 29998.8
$ws.Range("N50").Value = -31248.8
# Row 59
$ws.Range("H59").Value = 32479.4
$ws.Range("J59").Value = 34421.555
$ws.Range("L59").Value = 34421.555
$ws.Range("N59").Value = -36711.555
# Row 62
$ws.Range("H62").Value = 2906
$ws.Range("I62").Value = 2550
$ws.Range("J62").Value = 3024.6667
$ws.Range("K62").Value = 2550
$ws.Range("L62").Value = 3024.6667
$ws.Range("M62").Value = -1926
$ws.Range("N62").Value = -4272.6667
# Row 65
$ws.Range("H65").Value = 2906
$ws.Range("I65").Value = 2550
$ws.Range("J65").Value = 3024.6667
$ws.Range("K65").Value = 12750
$ws.Range("L65").Value = 15123.3335
$ws.Range("M65").Value = -9630
$ws.Range("N65").Value = -21363.3335
# Row 68
$ws.Range("H68").Value = 39999.11
$ws.Range("J68").Value = 39999.11
$ws.Range("L68").Value = 39999.11
$ws.Range("N68").Value = -41497.11
# Row 71
$ws.Range("H71").Value = 39999.11
$ws.Range("J71").Value = 39999.11
$ws.Range("L71").Value = 119997.33
$ws.Range("N71").Value = -127485.33
# Row 74
$ws.Range("H74").Value = 40000
$ws.Range("J74").Value = 40000
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41748
# Row 77
$ws.Range("H77").Value = 40000
$ws.Range("J77").Value = 40000
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -128736
# Row 92
$ws.Range("H92").Value = 41600
$ws.Range("I92").Value = 20000
$ws.Range("J92").Value = 47000
$ws.Range("K92").Value = 20000
$ws.Range("L92").Value = 47000
$ws.Range("M92").Value = -17504
$ws.Range("N92").Value = -51992
# Row 99
$ws.Range("H99").Value = 4870.8
$ws.Range("I99").Value = 3399.5
$ws.Range("K99").Value = 3399.5
$ws.Range("M99").Value = -1901.5
# Row 109
$ws.Range("H109").Value = 60118.5
$ws.Range("J109").Value = 60118.5
$ws.Range("L109").Value = 60118.5
$ws.Range("N109").Value = -62198.5
# Row 126
$ws.Range("H126").Value = 4870.8
$ws.Range("I126").Value = 3399.5
$ws.Range("K126").Value = 10198.5
$ws.Range("M126").Value = -7728.5

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 104.4
$ws.Range("I2").Value = 116.5
$ws.Range("J2").Value = 96.333336
$ws.Range("K2").Value = 699
$ws.Range("L2").Value = 578.000016
$ws.Range("M2").Value = -586
$ws.Range("N2").Value = -804.000016

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 63560.6
$ws.Range("I3").Value = 152000
$ws.Range("J3").Value = 4601
$ws.Range("K3").Value = 152000
$ws.Range("L3").Value = 4601
$ws.Range("M3").Value = -151884
$ws.Range("N3").Value = -4833
# Row 20
$ws.Range("H20").Value = 36000
$ws.Range("J20").Value = 36000
$ws.Range("L20").Value = 36000
$ws.Range("N20").Value = -36490
# Row 21
$ws.Range("H21").Value = 20669
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
# Row 30
$ws.Range("H30").Value = 20669
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
# Row 33
$ws.Range("H33").Value = 11499.5
$ws.Range("I33").Value = 10000
$ws.Range("K33").Value = 10000
$ws.Range("M33").Value = -9748
# Row 63
$ws.Range("H63").Value = 31000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 31000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 31000
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -32372
# Row 66
$ws.Range("H66").Value = 31000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 31000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 93000
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -99864
# Row 80
$ws.Range("H80").Value = 3199
$ws.Range("I80").Value = 2666.6667
$ws.Range("K80").Value = 2666.6667
$ws.Range("M80").Value = -1668.6667
# Row 83
$ws.Range("H83").Value = 3199
$ws.Range("I83").Value = 2666.6667
$ws.Range("K83").Value = 13333.3335
$ws.Range("M83").Value = -8341.333500000001
# Row 92
$ws.Range("H92").Value = 3135
$ws.Range("J92").Value = 3135
$ws.Range("L92").Value = 3135
$ws.Range("N92").Value = -6879
# Row 126
$ws.Range("H126").Value = 9253.272000000001
$ws.Range("I126").Value = 9253.272000000001
$ws.Range("K126").Value = 27759.816
$ws.Range("M126").Value = -25289.816

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 1662.875
$ws.Range("I4").Value = 1002
$ws.Range("K4").Value = 1002
$ws.Range("M4").Value = -889
# Row 51
$ws.Range("H51").Value = 4300
$ws.Range("I51").Value = 4300
$ws.Range("K51").Value = 4300
$ws.Range("M51").Value = -3790

Write-Output "done"